$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2026/12/29", "火", 13, 201),
    @("2026/12/29", "火", 16, 201),
    @("2026/12/29", "火", 19, 201),
    @("2026/12/29", "火", 23, 201),
    @("2026/12/30", "水", 2, 201),
    @("2026/12/30", "水", 5, 201),
    @("2026/12/30", "水", 8, 201),
    @("2026/12/30", "水", 13, 201),
    @("2026/12/30", "水", 16, 201),
    @("2026/12/30", "水", 22, 201),
    @("2026/12/31", "木", 2, 201),
    @("2026/12/31", "木", 6, 201),
    @("2026/12/31", "木", 10, 201),
    @("2026/12/31", "木", 12, 201),
    @("2026/12/31", "木", 14, 201),
    @("2026/12/31", "木", 22, 201),
    @("2027/01/01", "金", 2, 201),
    @("2027/01/01", "金", 5, 201),
    @("2027/01/01", "金", 13, 201),
    @("2027/01/01", "金", 16, 201),
    @("2027/01/01", "金", 19, 201),
    @("2027/01/02", "土", 1, 201),
    @("2027/01/02", "土", 5, 201),
    @("2027/01/02", "土", 8, 201),
    @("2027/01/02", "土", 13, 201),
    @("2027/01/02", "土", 16, 201),
    @("2027/01/02", "土", 19, 201),
    @("2027/01/02", "土", 22, 201),
    @("2027/01/03", "日", 1, 201),
    @("2027/01/03", "日", 4, 201),
    @("2027/01/03", "日", 7, 201),
    @("2027/01/03", "日", 13, 201),
    @("2027/01/03", "日", 16, 201),
    @("2027/01/03", "日", 19, 201),
    @("2027/01/03", "日", 22, 201),
    @("2027/01/04", "月", 2, 201),
    @("2027/01/04", "月", 4, 201),
    @("2027/01/04", "月", 7, 201),
    @("2027/01/04", "月", 13, 201),
    @("2027/01/04", "月", 22, 201),
    @("2027/01/05", "火", 2, 201),
    @("2027/01/05", "火", 7, 201)
)

$startRow = 571
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row[0]
    $cellA.Style = "Normal"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
